# Commit: "[update 2] minor changes"
# Updates the cached A (x) / B (y) value tables on both data sheets (100 rows each).
# Sheet 1 is named phi (unicode U+03C6); sheet 2 is named "z". Reference by index
# to avoid any non-ASCII literal/encoding issues in this script.
$wb = $excel.ActiveWorkbook
$wsPhi = $wb.Worksheets.Item(1)
$wsZ = $wb.Worksheets.Item(2)

# Sheet "phi": updated A (x) and B (y) columns, rows 1-100
$wsPhi.Cells.Item(1, 1).Value2 = 0
$wsPhi.Cells.Item(1, 2).Value2 = 0.104012573964497
$wsPhi.Cells.Item(2, 1).Value2 = 0.000005999999999999999
$wsPhi.Cells.Item(2, 2).Value2 = 0.1040117418639053
$wsPhi.Cells.Item(3, 1).Value2 = 0.000012
$wsPhi.Cells.Item(3, 2).Value2 = 0.1040059171597633
$wsPhi.Cells.Item(4, 1).Value2 = 0.000018
$wsPhi.Cells.Item(4, 2).Value2 = 0.1039901072485207
$wsPhi.Cells.Item(5, 1).Value2 = 0.000024
$wsPhi.Cells.Item(5, 2).Value2 = 0.1039593195266272
$wsPhi.Cells.Item(6, 1).Value2 = 0.00003
$wsPhi.Cells.Item(6, 2).Value2 = 0.1039085613905325
$wsPhi.Cells.Item(7, 1).Value2 = 0.00003599999999999999
$wsPhi.Cells.Item(7, 2).Value2 = 0.1038328402366864
$wsPhi.Cells.Item(8, 1).Value2 = 0.00004199999999999999
$wsPhi.Cells.Item(8, 2).Value2 = 0.1037271634615384
$wsPhi.Cells.Item(9, 1).Value2 = 0.00004799999999999999
$wsPhi.Cells.Item(9, 2).Value2 = 0.1035865384615385
$wsPhi.Cells.Item(10, 1).Value2 = 0.00005399999999999998
$wsPhi.Cells.Item(10, 2).Value2 = 0.1034059726331361
$wsPhi.Cells.Item(11, 1).Value2 = 0.00005999999999999998
$wsPhi.Cells.Item(11, 2).Value2 = 0.1031804733727811
$wsPhi.Cells.Item(12, 1).Value2 = 0.00006599999999999998
$wsPhi.Cells.Item(12, 2).Value2 = 0.1029050480769231
$wsPhi.Cells.Item(13, 1).Value2 = 0.00007199999999999997
$wsPhi.Cells.Item(13, 2).Value2 = 0.1025747041420118
$wsPhi.Cells.Item(14, 1).Value2 = 0.00007799999999999997
$wsPhi.Cells.Item(14, 2).Value2 = 0.102184448964497
$wsPhi.Cells.Item(15, 1).Value2 = 0.00008399999999999997
$wsPhi.Cells.Item(15, 2).Value2 = 0.1017292899408284
$wsPhi.Cells.Item(16, 1).Value2 = 0.00008999999999999997
$wsPhi.Cells.Item(16, 2).Value2 = 0.1012042344674556
$wsPhi.Cells.Item(17, 1).Value2 = 0.00009599999999999996
$wsPhi.Cells.Item(17, 2).Value2 = 0.1006042899408284
$wsPhi.Cells.Item(18, 1).Value2 = 0.000102
$wsPhi.Cells.Item(18, 2).Value2 = 0.09992446375739644
$wsPhi.Cells.Item(19, 1).Value2 = 0.000108
$wsPhi.Cells.Item(19, 2).Value2 = 0.09915976331360946
$wsPhi.Cells.Item(20, 1).Value2 = 0.000114
$wsPhi.Cells.Item(20, 2).Value2 = 0.09830519600591717
$wsPhi.Cells.Item(21, 1).Value2 = 0.0001199999999999999
$wsPhi.Cells.Item(21, 2).Value2 = 0.09735576923076922
$wsPhi.Cells.Item(22, 1).Value2 = 0.0001259999999999999
$wsPhi.Cells.Item(22, 2).Value2 = 0.09630649038461539
$wsPhi.Cells.Item(23, 1).Value2 = 0.000132
$wsPhi.Cells.Item(23, 2).Value2 = 0.09515236686390531
$wsPhi.Cells.Item(24, 1).Value2 = 0.000138
$wsPhi.Cells.Item(24, 2).Value2 = 0.09388840606508876
$wsPhi.Cells.Item(25, 1).Value2 = 0.000144
$wsPhi.Cells.Item(25, 2).Value2 = 0.09250961538461538
$wsPhi.Cells.Item(26, 1).Value2 = 0.00015
$wsPhi.Cells.Item(26, 2).Value2 = 0.0910110022189349
$wsPhi.Cells.Item(27, 1).Value2 = 0.000156
$wsPhi.Cells.Item(27, 2).Value2 = 0.08938757396449704
$wsPhi.Cells.Item(28, 1).Value2 = 0.000162
$wsPhi.Cells.Item(28, 2).Value2 = 0.08763433801775146
$wsPhi.Cells.Item(29, 1).Value2 = 0.000168
$wsPhi.Cells.Item(29, 2).Value2 = 0.0857463017751479
$wsPhi.Cells.Item(30, 1).Value2 = 0.000174
$wsPhi.Cells.Item(30, 2).Value2 = 0.08371847263313607
$wsPhi.Cells.Item(31, 1).Value2 = 0.00018
$wsPhi.Cells.Item(31, 2).Value2 = 0.08154585798816566
$wsPhi.Cells.Item(32, 1).Value2 = 0.000186
$wsPhi.Cells.Item(32, 2).Value2 = 0.07922346523668636
$wsPhi.Cells.Item(33, 1).Value2 = 0.0001920000000000001
$wsPhi.Cells.Item(33, 2).Value2 = 0.07674630177514788
$wsPhi.Cells.Item(34, 1).Value2 = 0.0001980000000000001
$wsPhi.Cells.Item(34, 2).Value2 = 0.07410937499999995
$wsPhi.Cells.Item(35, 1).Value2 = 0.0002040000000000001
$wsPhi.Cells.Item(35, 2).Value2 = 0.07130769230769225
$wsPhi.Cells.Item(36, 1).Value2 = 0.0002100000000000001
$wsPhi.Cells.Item(36, 2).Value2 = 0.06833626109467449
$wsPhi.Cells.Item(37, 1).Value2 = 0.0002160000000000001
$wsPhi.Cells.Item(37, 2).Value2 = 0.06519008875739637
$wsPhi.Cells.Item(38, 1).Value2 = 0.0002220000000000001
$wsPhi.Cells.Item(38, 2).Value2 = 0.06186418269230761
$wsPhi.Cells.Item(39, 1).Value2 = 0.0002280000000000001
$wsPhi.Cells.Item(39, 2).Value2 = 0.0583535502958579
$wsPhi.Cells.Item(40, 1).Value2 = 0.0002340000000000001
$wsPhi.Cells.Item(40, 2).Value2 = 0.05465319896449694
$wsPhi.Cells.Item(41, 1).Value2 = 0.0002400000000000001
$wsPhi.Cells.Item(41, 2).Value2 = 0.05075813609467444
$wsPhi.Cells.Item(42, 1).Value2 = 0.0002460000000000001
$wsPhi.Cells.Item(42, 2).Value2 = 0.04666336908284014
$wsPhi.Cells.Item(43, 1).Value2 = 0.0002520000000000001
$wsPhi.Cells.Item(43, 2).Value2 = 0.04236390532544369
$wsPhi.Cells.Item(44, 1).Value2 = 0.0002580000000000001
$wsPhi.Cells.Item(44, 2).Value2 = 0.03785475221893483
$wsPhi.Cells.Item(45, 1).Value2 = 0.0002640000000000001
$wsPhi.Cells.Item(45, 2).Value2 = 0.03313091715976323
$wsPhi.Cells.Item(46, 1).Value2 = 0.0002700000000000001
$wsPhi.Cells.Item(46, 2).Value2 = 0.02818740754437863
$wsPhi.Cells.Item(47, 1).Value2 = 0.000276
$wsPhi.Cells.Item(47, 2).Value2 = 0.02301923076923071
$wsPhi.Cells.Item(48, 1).Value2 = 0.000282
$wsPhi.Cells.Item(48, 2).Value2 = 0.01762139423076918
$wsPhi.Cells.Item(49, 1).Value2 = 0.000288
$wsPhi.Cells.Item(49, 2).Value2 = 0.01198890532544376
$wsPhi.Cells.Item(50, 1).Value2 = 0.000294
$wsPhi.Cells.Item(50, 2).Value2 = 0.006116771449704132
$wsPhi.Cells.Item(51, 1).Value2 = 0.0003
$wsPhi.Cells.Item(51, 2).Value2 = 0
$wsPhi.Cells.Item(52, 1).Value2 = 0.000306
$wsPhi.Cells.Item(52, 2).Value2 = -0.006366401627218911
$wsPhi.Cells.Item(53, 1).Value2 = 0.0003119999999999999
$wsPhi.Cells.Item(53, 2).Value2 = -0.01298742603550293
$wsPhi.Cells.Item(54, 1).Value2 = 0.0003179999999999999
$wsPhi.Cells.Item(54, 2).Value2 = -0.01986806582840229
$wsPhi.Cells.Item(55, 1).Value2 = 0.0003239999999999999
$wsPhi.Cells.Item(55, 2).Value2 = -0.02701331360946738
$wsPhi.Cells.Item(56, 1).Value2 = 0.0003299999999999999
$wsPhi.Cells.Item(56, 2).Value2 = -0.0344281619822484
$wsPhi.Cells.Item(57, 1).Value2 = 0.0003359999999999999
$wsPhi.Cells.Item(57, 2).Value2 = -0.04211760355029573
$wsPhi.Cells.Item(58, 1).Value2 = 0.0003419999999999999
$wsPhi.Cells.Item(58, 2).Value2 = -0.05008663091715961
$wsPhi.Cells.Item(59, 1).Value2 = 0.0003479999999999998
$wsPhi.Cells.Item(59, 2).Value2 = -0.05834023668639034
$wsPhi.Cells.Item(60, 1).Value2 = 0.0003539999999999998
$wsPhi.Cells.Item(60, 2).Value2 = -0.06688341346153824
$wsPhi.Cells.Item(61, 1).Value2 = 0.0003599999999999998
$wsPhi.Cells.Item(61, 2).Value2 = -0.07572115384615358
$wsPhi.Cells.Item(62, 1).Value2 = 0.0003659999999999998
$wsPhi.Cells.Item(62, 2).Value2 = -0.08485845044378668
$wsPhi.Cells.Item(63, 1).Value2 = 0.0003719999999999998
$wsPhi.Cells.Item(63, 2).Value2 = -0.09430029585798783
$wsPhi.Cells.Item(64, 1).Value2 = 0.0003779999999999998
$wsPhi.Cells.Item(64, 2).Value2 = -0.1040516826923073
$wsPhi.Cells.Item(65, 1).Value2 = 0.0003839999999999997
$wsPhi.Cells.Item(65, 2).Value2 = -0.1141176035502954
$wsPhi.Cells.Item(66, 1).Value2 = 0.0003899999999999997
$wsPhi.Cells.Item(66, 2).Value2 = -0.1245030510355025
$wsPhi.Cells.Item(67, 1).Value2 = 0.0003959999999999997
$wsPhi.Cells.Item(67, 2).Value2 = -0.1352130177514788
$wsPhi.Cells.Item(68, 1).Value2 = 0.0004019999999999997
$wsPhi.Cells.Item(68, 2).Value2 = -0.1462524963017746
$wsPhi.Cells.Item(69, 1).Value2 = 0.0004079999999999997
$wsPhi.Cells.Item(69, 2).Value2 = -0.1576264792899403
$wsPhi.Cells.Item(70, 1).Value2 = 0.0004139999999999997
$wsPhi.Cells.Item(70, 2).Value2 = -0.169339959319526
$wsPhi.Cells.Item(71, 1).Value2 = 0.0004199999999999996
$wsPhi.Cells.Item(71, 2).Value2 = -0.1813979289940821
$wsPhi.Cells.Item(72, 1).Value2 = 0.0004259999999999996
$wsPhi.Cells.Item(72, 2).Value2 = -0.193805380917159
$wsPhi.Cells.Item(73, 1).Value2 = 0.0004319999999999996
$wsPhi.Cells.Item(73, 2).Value2 = -0.2065673076923069
$wsPhi.Cells.Item(74, 1).Value2 = 0.0004379999999999996
$wsPhi.Cells.Item(74, 2).Value2 = -0.219688701923076
$wsPhi.Cells.Item(75, 1).Value2 = 0.0004439999999999996
$wsPhi.Cells.Item(75, 2).Value2 = -0.2331745562130168
$wsPhi.Cells.Item(76, 1).Value2 = 0.0004499999999999996
$wsPhi.Cells.Item(76, 2).Value2 = -0.2470298631656794
$wsPhi.Cells.Item(77, 1).Value2 = 0.0004559999999999995
$wsPhi.Cells.Item(77, 2).Value2 = -0.2612596153846143
$wsPhi.Cells.Item(78, 1).Value2 = 0.0004619999999999995
$wsPhi.Cells.Item(78, 2).Value2 = -0.2758688054733716
$wsPhi.Cells.Item(79, 1).Value2 = 0.0004679999999999995
$wsPhi.Cells.Item(79, 2).Value2 = -0.2908624260355017
$wsPhi.Cells.Item(80, 1).Value2 = 0.0004739999999999995
$wsPhi.Cells.Item(80, 2).Value2 = -0.3062454696745549
$wsPhi.Cells.Item(81, 1).Value2 = 0.0004799999999999995
$wsPhi.Cells.Item(81, 2).Value2 = -0.3220229289940815
$wsPhi.Cells.Item(82, 1).Value2 = 0.0004859999999999995
$wsPhi.Cells.Item(82, 2).Value2 = -0.3381997965976317
$wsPhi.Cells.Item(83, 1).Value2 = 0.0004919999999999995
$wsPhi.Cells.Item(83, 2).Value2 = -0.354781065088756
$wsPhi.Cells.Item(84, 1).Value2 = 0.0004979999999999995
$wsPhi.Cells.Item(84, 2).Value2 = -0.3717717270710046
$wsPhi.Cells.Item(85, 1).Value2 = 0.0005039999999999996
$wsPhi.Cells.Item(85, 2).Value2 = -0.3891767751479278
$wsPhi.Cells.Item(86, 1).Value2 = 0.0005099999999999996
$wsPhi.Cells.Item(86, 2).Value2 = -0.4070012019230758
$wsPhi.Cells.Item(87, 1).Value2 = 0.0005159999999999996
$wsPhi.Cells.Item(87, 2).Value2 = -0.4252499999999989
$wsPhi.Cells.Item(88, 1).Value2 = 0.0005219999999999997
$wsPhi.Cells.Item(88, 2).Value2 = -0.4439281619822475
$wsPhi.Cells.Item(89, 1).Value2 = 0.0005279999999999997
$wsPhi.Cells.Item(89, 2).Value2 = -0.463040680473372
$wsPhi.Cells.Item(90, 1).Value2 = 0.0005339999999999998
$wsPhi.Cells.Item(90, 2).Value2 = -0.4825925480769224
$wsPhi.Cells.Item(91, 1).Value2 = 0.0005399999999999998
$wsPhi.Cells.Item(91, 2).Value2 = -0.5025887573964491
$wsPhi.Cells.Item(92, 1).Value2 = 0.0005459999999999998
$wsPhi.Cells.Item(92, 2).Value2 = -0.5230343010355024
$wsPhi.Cells.Item(93, 1).Value2 = 0.0005519999999999999
$wsPhi.Cells.Item(93, 2).Value2 = -0.5439341715976328
$wsPhi.Cells.Item(94, 1).Value2 = 0.0005579999999999999
$wsPhi.Cells.Item(94, 2).Value2 = -0.5652933616863902
$wsPhi.Cells.Item(95, 1).Value2 = 0.0005639999999999999
$wsPhi.Cells.Item(95, 2).Value2 = -0.5871168639053254
$wsPhi.Cells.Item(96, 1).Value2 = 0.00057
$wsPhi.Cells.Item(96, 2).Value2 = -0.6094096708579883
$wsPhi.Cells.Item(97, 1).Value2 = 0.000576
$wsPhi.Cells.Item(97, 2).Value2 = -0.6321767751479291
$wsPhi.Cells.Item(98, 1).Value2 = 0.000582
$wsPhi.Cells.Item(98, 2).Value2 = -0.6554231693786986
$wsPhi.Cells.Item(99, 1).Value2 = 0.0005880000000000001
$wsPhi.Cells.Item(99, 2).Value2 = -0.6791538461538467
$wsPhi.Cells.Item(100, 1).Value2 = 0.0005940000000000001
$wsPhi.Cells.Item(100, 2).Value2 = -0.7033737980769237

# Sheet "z": updated A (x) and B (y) columns, rows 1-100
$wsZ.Cells.Item(1, 1).Value2 = 0
$wsZ.Cells.Item(1, 2).Value2 = -0.00002340282914201183
$wsZ.Cells.Item(2, 1).Value2 = 0.000005999999999999999
$wsZ.Cells.Item(2, 2).Value2 = -0.00002277875494637573
$wsZ.Cells.Item(3, 1).Value2 = 0.000012
$wsZ.Cells.Item(3, 2).Value2 = -0.00002215469822485206
$wsZ.Cells.Item(4, 1).Value2 = 0.000018
$wsZ.Cells.Item(4, 2).Value2 = -0.00002153070391087277
$wsZ.Cells.Item(5, 1).Value2 = 0.000024
$wsZ.Cells.Item(5, 2).Value2 = -0.00002090684689349112
$wsZ.Cells.Item(6, 1).Value2 = 0.00003
$wsZ.Cells.Item(6, 2).Value2 = -0.00002028323201738165
$wsZ.Cells.Item(7, 1).Value2 = 0.00003599999999999999
$wsZ.Cells.Item(7, 2).Value2 = -0.00001965999408284023
$wsZ.Cells.Item(8, 1).Value2 = 0.00004199999999999999
$wsZ.Cells.Item(8, 2).Value2 = -0.00001903729784578401
$wsZ.Cells.Item(9, 1).Value2 = 0.00004799999999999999
$wsZ.Cells.Item(9, 2).Value2 = -0.00001841533801775147
$wsZ.Cells.Item(10, 1).Value2 = 0.00005399999999999998
$wsZ.Cells.Item(10, 2).Value2 = -0.00001779433926590236
$wsZ.Cells.Item(11, 1).Value2 = 0.00005999999999999998
$wsZ.Cells.Item(11, 2).Value2 = -0.00001717455621301775
$wsZ.Cells.Item(12, 1).Value2 = 0.00006599999999999998
$wsZ.Cells.Item(12, 2).Value2 = -0.00001655627343749999
$wsZ.Cells.Item(13, 1).Value2 = 0.00007199999999999997
$wsZ.Cells.Item(13, 2).Value2 = -0.00001593980547337278
$wsZ.Cells.Item(14, 1).Value2 = 0.00007799999999999997
$wsZ.Cells.Item(14, 2).Value2 = -0.00001532549681028106
$wsZ.Cells.Item(15, 1).Value2 = 0.00008399999999999997
$wsZ.Cells.Item(15, 2).Value2 = -0.00001471372189349112
$wsZ.Cells.Item(16, 1).Value2 = 0.00008999999999999997
$wsZ.Cells.Item(16, 2).Value2 = -0.00001410488512389053
$wsZ.Cells.Item(17, 1).Value2 = 0.00009599999999999996
$wsZ.Cells.Item(17, 2).Value2 = -0.00001349942085798816
$wsZ.Cells.Item(18, 1).Value2 = 0.000102
$wsZ.Cells.Item(18, 2).Value2 = -0.0000128977934079142
$wsZ.Cells.Item(19, 1).Value2 = 0.000108
$wsZ.Cells.Item(19, 2).Value2 = -0.00001230049704142011
$wsZ.Cells.Item(20, 1).Value2 = 0.000114
$wsZ.Cells.Item(20, 2).Value2 = -0.0000117080559818787
$wsZ.Cells.Item(21, 1).Value2 = 0.0001199999999999999
$wsZ.Cells.Item(21, 2).Value2 = -0.00001112102440828402
$wsZ.Cells.Item(22, 1).Value2 = 0.0001259999999999999
$wsZ.Cells.Item(22, 2).Value2 = -0.00001053998645525148
$wsZ.Cells.Item(23, 1).Value2 = 0.000132
$wsZ.Cells.Item(23, 2).Value2 = -0.000009965556213017749
$wsZ.Cells.Item(24, 1).Value2 = 0.000138
$wsZ.Cells.Item(24, 2).Value2 = -0.000009398377727440823
$wsZ.Cells.Item(25, 1).Value2 = 0.000144
$wsZ.Cells.Item(25, 2).Value2 = -0.000008839124999999995
$wsZ.Cells.Item(26, 1).Value2 = 0.00015
$wsZ.Cells.Item(26, 2).Value2 = -0.000008288501987795854
$wsZ.Cells.Item(27, 1).Value2 = 0.000156
$wsZ.Cells.Item(27, 2).Value2 = -0.000007747242603550288
$wsZ.Cells.Item(28, 1).Value2 = 0.000162
$wsZ.Cells.Item(28, 2).Value2 = -0.000007216110715606505
$wsZ.Cells.Item(29, 1).Value2 = 0.000168
$wsZ.Cells.Item(29, 2).Value2 = -0.000006695900147928985
$wsZ.Cells.Item(30, 1).Value2 = 0.000174
$wsZ.Cells.Item(30, 2).Value2 = -0.000006187434680103541
$wsZ.Cells.Item(31, 1).Value2 = 0.00018
$wsZ.Cells.Item(31, 2).Value2 = -0.00000569156804733727
$wsZ.Cells.Item(32, 1).Value2 = 0.000186
$wsZ.Cells.Item(32, 2).Value2 = -0.000005209183940458572
$wsZ.Cells.Item(33, 1).Value2 = 0.0001920000000000001
$wsZ.Cells.Item(33, 2).Value2 = -0.000004741196005917148
$wsZ.Cells.Item(34, 1).Value2 = 0.0001980000000000001
$wsZ.Cells.Item(34, 2).Value2 = -0.000004288547845784015
$wsZ.Cells.Item(35, 1).Value2 = 0.0002040000000000001
$wsZ.Cells.Item(35, 2).Value2 = -0.000003852213017751466
$wsZ.Cells.Item(36, 1).Value2 = 0.0002100000000000001
$wsZ.Cells.Item(36, 2).Value2 = -0.000003433195035133126
$wsZ.Cells.Item(37, 1).Value2 = 0.0002160000000000001
$wsZ.Cells.Item(37, 2).Value2 = -0.000003032527366863892
$wsZ.Cells.Item(38, 1).Value2 = 0.0002220000000000001
$wsZ.Cells.Item(38, 2).Value2 = -0.000002651273437499987
$wsZ.Cells.Item(39, 1).Value2 = 0.0002280000000000001
$wsZ.Cells.Item(39, 2).Value2 = -0.00000229052662721892
$wsZ.Cells.Item(40, 1).Value2 = 0.0002340000000000001
$wsZ.Cells.Item(40, 2).Value2 = -0.000001951410271819514
$wsZ.Cells.Item(41, 1).Value2 = 0.0002400000000000001
$wsZ.Cells.Item(41, 2).Value2 = -0.000001635077662721879
$wsZ.Cells.Item(42, 1).Value2 = 0.0002460000000000001
$wsZ.Cells.Item(42, 2).Value2 = -0.000001342712046967447
$wsZ.Cells.Item(43, 1).Value2 = 0.0002520000000000001
$wsZ.Cells.Item(43, 2).Value2 = -0.000001075526627218924
$wsZ.Cells.Item(44, 1).Value2 = 0.0002580000000000001
$wsZ.Cells.Item(44, 2).Value2 = -0.0000008347645617603467
$wsZ.Cells.Item(45, 1).Value2 = 0.0002640000000000001
$wsZ.Cells.Item(45, 2).Value2 = -0.0000006216989644970331
$wsZ.Cells.Item(46, 1).Value2 = 0.0002700000000000001
$wsZ.Cells.Item(46, 2).Value2 = -0.000000437632904955617
$wsZ.Cells.Item(47, 1).Value2 = 0.000276
$wsZ.Cells.Item(47, 2).Value2 = -0.0000002838994082840209
$wsZ.Cells.Item(48, 1).Value2 = 0.000282
$wsZ.Cells.Item(48, 2).Value2 = -0.0000001618614552514737
$wsZ.Cells.Item(49, 1).Value2 = 0.000288
$wsZ.Cells.Item(49, 2).Value2 = -0.00000007291198224851449
$wsZ.Cells.Item(50, 1).Value2 = 0.000294
$wsZ.Cells.Item(50, 2).Value2 = -0.00000001847388128697958
$wsZ.Cells.Item(51, 1).Value2 = 0.0003
$wsZ.Cells.Item(51, 2).Value2 = 0
$wsZ.Cells.Item(52, 1).Value2 = 0.000306
$wsZ.Cells.Item(52, 2).Value2 = -0.00000001897314164200437
$wsZ.Cells.Item(53, 1).Value2 = 0.0003119999999999999
$wsZ.Cells.Item(53, 2).Value2 = -0.00000007690606508875231
$wsZ.Cells.Item(54, 1).Value2 = 0.0003179999999999999
$wsZ.Cells.Item(54, 2).Value2 = -0.0000001753414848372768
$wsZ.Cells.Item(55, 1).Value2 = 0.0003239999999999999
$wsZ.Cells.Item(55, 2).Value2 = -0.0000003158520710059082
$wsZ.Cells.Item(56, 1).Value2 = 0.0003299999999999999
$wsZ.Cells.Item(56, 2).Value2 = -0.0000005000404493343113
$wsZ.Cells.Item(57, 1).Value2 = 0.0003359999999999999
$wsZ.Cells.Item(57, 2).Value2 = -0.0000007295392011834237
$wsZ.Cells.Item(58, 1).Value2 = 0.0003419999999999999
$wsZ.Cells.Item(58, 2).Value2 = -0.000001006010863535493
$wsZ.Cells.Item(59, 1).Value2 = 0.0003479999999999998
$wsZ.Cells.Item(59, 2).Value2 = -0.000001331147928994068
$wsZ.Cells.Item(60, 1).Value2 = 0.0003539999999999998
$wsZ.Cells.Item(60, 2).Value2 = -0.00000170667284578401
$wsZ.Cells.Item(61, 1).Value2 = 0.0003599999999999998
$wsZ.Cells.Item(61, 2).Value2 = -0.000002134338017751464
$wsZ.Cells.Item(62, 1).Value2 = 0.0003659999999999998
$wsZ.Cells.Item(62, 2).Value2 = -0.000002615925804363887
$wsZ.Cells.Item(63, 1).Value2 = 0.0003719999999999998
$wsZ.Cells.Item(63, 2).Value2 = -0.000003153248520710033
$wsZ.Cells.Item(64, 1).Value2 = 0.0003779999999999998
$wsZ.Cells.Item(64, 2).Value2 = -0.000003748148437499972
$wsZ.Cells.Item(65, 1).Value2 = 0.0003839999999999997
$wsZ.Cells.Item(65, 2).Value2 = -0.000004402497781065059
$wsZ.Cells.Item(66, 1).Value2 = 0.0003899999999999997
$wsZ.Cells.Item(66, 2).Value2 = -0.000005118198733357949
$wsZ.Cells.Item(67, 1).Value2 = 0.0003959999999999997
$wsZ.Cells.Item(67, 2).Value2 = -0.00000589718343195262
$wsZ.Cells.Item(68, 1).Value2 = 0.0004019999999999997
$wsZ.Cells.Item(68, 2).Value2 = -0.000006741413970044332
$wsZ.Cells.Item(69, 1).Value2 = 0.0004079999999999997
$wsZ.Cells.Item(69, 2).Value2 = -0.000007652882396449653
$wsZ.Cells.Item(70, 1).Value2 = 0.0004139999999999997
$wsZ.Cells.Item(70, 2).Value2 = -0.000008633610715606448
$wsZ.Cells.Item(71, 1).Value2 = 0.0004199999999999996
$wsZ.Cells.Item(71, 2).Value2 = -0.0000096856508875739
$wsZ.Cells.Item(72, 1).Value2 = 0.0004259999999999996
$wsZ.Cells.Item(72, 2).Value2 = -0.00001081108482803247
$wsZ.Cells.Item(73, 1).Value2 = 0.0004319999999999996
$wsZ.Cells.Item(73, 2).Value2 = -0.00001201202440828394
$wsZ.Cells.Item(74, 1).Value2 = 0.0004379999999999996
$wsZ.Cells.Item(74, 2).Value2 = -0.00001329061145525139
$wsZ.Cells.Item(75, 1).Value2 = 0.0004439999999999996
$wsZ.Cells.Item(75, 2).Value2 = -0.00001464901775147919
$wsZ.Cells.Item(76, 1).Value2 = 0.0004499999999999996
$wsZ.Cells.Item(76, 2).Value2 = -0.00001608944503513303
$wsZ.Cells.Item(77, 1).Value2 = 0.0004559999999999995
$wsZ.Cells.Item(77, 2).Value2 = -0.00001761412499999988
$wsZ.Cells.Item(78, 1).Value2 = 0.0004619999999999995
$wsZ.Cells.Item(78, 2).Value2 = -0.00001922531929548803
$wsZ.Cells.Item(79, 1).Value2 = 0.0004679999999999995
$wsZ.Cells.Item(79, 2).Value2 = -0.00002092531952662708
$wsZ.Cells.Item(80, 1).Value2 = 0.0004739999999999995
$wsZ.Cells.Item(80, 2).Value2 = -0.00002271644725406789
$wsZ.Cells.Item(81, 1).Value2 = 0.0004799999999999995
$wsZ.Cells.Item(81, 2).Value2 = -0.00002460105399408267
$wsZ.Cells.Item(82, 1).Value2 = 0.0004859999999999995
$wsZ.Cells.Item(82, 2).Value2 = -0.00002658152121856491
$wsZ.Cells.Item(83, 1).Value2 = 0.0004919999999999995
$wsZ.Cells.Item(83, 2).Value2 = -0.00002866026035502941
$wsZ.Cells.Item(84, 1).Value2 = 0.0004979999999999995
$wsZ.Cells.Item(84, 2).Value2 = -0.00003083971278661225
$wsZ.Cells.Item(85, 1).Value2 = 0.0005039999999999996
$wsZ.Cells.Item(85, 2).Value2 = -0.00003312234985207084
$wsZ.Cells.Item(86, 1).Value2 = 0.0005099999999999996
$wsZ.Cells.Item(86, 2).Value2 = -0.00003551067284578386
$wsZ.Cells.Item(87, 1).Value2 = 0.0005159999999999996
$wsZ.Cells.Item(87, 2).Value2 = -0.00003800721301775133
$wsZ.Cells.Item(88, 1).Value2 = 0.0005219999999999997
$wsZ.Cells.Item(88, 2).Value2 = -0.00004061453157359454
$wsZ.Cells.Item(89, 1).Value2 = 0.0005279999999999997
$wsZ.Cells.Item(89, 2).Value2 = -0.00004333521967455609
$wsZ.Cells.Item(90, 1).Value2 = 0.0005339999999999998
$wsZ.Cells.Item(90, 2).Value2 = -0.00004617189843749988
$wsZ.Cells.Item(91, 1).Value2 = 0.0005399999999999998
$wsZ.Cells.Item(91, 2).Value2 = -0.00004912721893491115
$wsZ.Cells.Item(92, 1).Value2 = 0.0005459999999999998
$wsZ.Cells.Item(92, 2).Value2 = -0.00005220386219489636
$wsZ.Cells.Item(93, 1).Value2 = 0.0005519999999999999
$wsZ.Cells.Item(93, 2).Value2 = -0.00005540453920118336
$wsZ.Cells.Item(94, 1).Value2 = 0.0005579999999999999
$wsZ.Cells.Item(94, 2).Value2 = -0.00005873199089312125
$wsZ.Cells.Item(95, 1).Value2 = 0.0005639999999999999
$wsZ.Cells.Item(95, 2).Value2 = -0.00006218898816568045
$wsZ.Cells.Item(96, 1).Value2 = 0.00057
$wsZ.Cells.Item(96, 2).Value2 = -0.00006577833186945265
$wsZ.Cells.Item(97, 1).Value2 = 0.000576
$wsZ.Cells.Item(97, 2).Value2 = -0.00006950285281065091
$wsZ.Cells.Item(98, 1).Value2 = 0.000582
$wsZ.Cells.Item(98, 2).Value2 = -0.00007336541175110951
$wsZ.Cells.Item(99, 1).Value2 = 0.0005880000000000001
$wsZ.Cells.Item(99, 2).Value2 = -0.00007736889940828408
$wsZ.Cells.Item(100, 1).Value2 = 0.0005940000000000001
$wsZ.Cells.Item(100, 2).Value2 = -0.00008151623645525156

Write-Host "Updated 200 rows across both data sheets"
